$d = $word.ActiveDocument
$wNs = 'xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"'

function Append-ParaXml($innerXml) {
    $r = $d.Content
    $r.Collapse(0)
    $xml = '<w:p ' + $wNs + '>' + $innerXml + '</w:p>'
    [void]$r.InsertXML($xml)
}

# Paragraph: separator line
Append-ParaXml('<w:r><w:t>--------------------------------------------------------------------------------------------------------------------------------------</w:t></w:r>')

# Paragraph 1
Append-ParaXml('<w:r><w:t xml:space="preserve">1. You can change the route to change the root </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>url</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> for the website.</w:t></w:r>')

# Paragraph 2
Append-ParaXml('<w:r><w:t>2.  app/views/layouts/</w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>application.html.erb</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> contains the page layout used for all views.</w:t></w:r>')

# Paragraph 3
Append-ParaXml('<w:r><w:lastRenderedPageBreak/><w:t>3.</w:t></w:r><w:r><w:t xml:space="preserve"> </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>Turbolinks</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> increases the speeds up how long it takes to switch between pages on your website by keeping the current page instance and only changing the body and title of the page, instead of recompiling the </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>javascript</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> and </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>css</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> every time.</w:t></w:r>')

# Paragraph 4
Append-ParaXml('<w:r><w:t xml:space="preserve">4. </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>csrf_meta_tags</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t>() helps to prevent cross site scripting attacks by using an authenticity token to verify that a request is valid.</w:t></w:r>')

# Paragraph 5
Append-ParaXml('<w:r><w:t xml:space="preserve">5. </w:t></w:r><w:r><w:t>yield inserts the actual page content into the layout page.</w:t></w:r>')

# Paragraph 6
Append-ParaXml('<w:r><w:t>6. You might want to test that certain parts of the page are formatted correctly, and that the page itself responds.</w:t></w:r><w:r><w:t xml:space="preserve"> </w:t></w:r>')

Write-Output "done"
